$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1W", 0.024),
    @("4C", 0.994),
    @("4I", 1.015666667),
    @("4W", 0.666),
    @("5C", 1.0095),
    @("5F", 1.0485),
    @("5I", 1.049),
    @("5L", 1.0235),
    @("5N", 1.0465),
    @("5V", 1.0095),
    @("5Y", 1.0605),
    @("25Q", 0.994)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
